# Adds the 4 new numpad "A"/"B"/"C"/"D" key shapes to slide 2 ("34-key
# numpad" layout), matching the shapes already present for the digit keys
# (e.g. the "7" key). The new keys sit to the right of "9", same row.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points (a
# 32-bit float) while the underlying OOXML stores EMU (1 pt = 12700 EMU).
# A tiny epsilon is added before the EMU->pt conversion to counteract
# float32 truncation when PowerPoint converts the point value back to EMU.
function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00002
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Use the existing "7" key (first digit-key rectangle) as the template: it
# already has the right fill (tx1), text color (bg1), font (Fira Code,
# 11pt) and style refs that the new keys need.
$template = $s.Shapes.Item(1)

# PowerPoint hands out shape ids from a per-slide, ever-increasing counter
# (skipping ids already used on the slide). Fast-forward that counter with
# throwaway duplicates so the first new key we keep gets id 87, continuing
# on from the deck's existing shapes (which end at id 86 / "Rectangle 85").
# (Guarded with a loop cap so a differently-numbered deck can't spin forever.)
$targetFirstId = 87
$current = $template.Duplicate().Item(1)
$guard = 0
while (($current.Id -lt $targetFirstId) -and ($guard -lt 10000)) {
    $current.Delete()
    $current = $template.Duplicate().Item(1)
    $guard++
}

$letters = @("A", "B", "C", "D")
$xOffsets = @(3710070, 4593990, 5477910, 6361830)
$yOffset = 3912286
$width = 576000
$height = 576000
$startName = 86

for ($i = 0; $i -lt $letters.Length; $i++) {
    if ($i -eq 0) {
        $newShape = $current
    } else {
        $newShape = $template.Duplicate().Item(1)
    }

    $newShape.Left = EmuToPt($xOffsets[$i])
    $newShape.Top = EmuToPt($yOffset)
    $newShape.Width = EmuToPt($width)
    $newShape.Height = EmuToPt($height)
    $newShape.Name = "Rectangle " + ($startName + $i)

    $newShape.TextFrame.TextRange.Text = $letters[$i]
}
